$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.932.68'
$ws.Range("E2").Value = '  +2.97%  '
$ws.Range("D3").Value = '3.640.01'
$ws.Range("E3").Value = '  +6.63%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''588.77'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '''181.14'
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").Value = '3.630.68'
$ws.Range("E7").Value = '  +6.57%  '
$ws.Range("E8").Value = '  +2.53%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("E11").Value = '  +2.36%  '
$ws.Range("D12").Value = '''49.90'
$ws.Range("E12").Value = '  +2.43%  '
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '''683.73'
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = '4.227.60'
$ws.Range("E15").Value = '  +6.62%  '
$ws.Range("D16").Value = '''9.05'
$ws.Range("E16").Value = '  +3.78%  '
$ws.Range("D17").Value = '3.660.11'
$ws.Range("E17").Value = '  +7.06%  '
$ws.Range("D18").Value = '71.969.71'
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '''18.34'
$ws.Range("E20").Value = '  +2.84%  '
$ws.Range("D21").Value = '''11.63'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").Value = '''0.940'
$ws.Range("E22").Value = '  +2.56%  '
$ws.Range("D23").Value = '''5.93'
$ws.Range("E23").Value = '  +10.21%  '
$ws.Range("D24").Value = '''17.85'
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("D25").Value = '''103.37'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("E27").Value = '  +4.85%  '
$ws.Range("D28").Value = '''10.02'
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("E29").Value = '  +2.87%  '
$ws.Range("D30").Value = '''9.26'
$ws.Range("E30").Value = '  +4.33%  '
$ws.Range("E31").Value = '  +5.07%  '
$ws.Range("E32").Value = '  +14.24%  '
$ws.Range("D33").Value = '''591.46'
$ws.Range("E33").Value = '  +5.98%  '
$ws.Range("D34").Value = '''11.35'
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("D36").Value = '''59.40'
$ws.Range("E36").Value = '  +1.07%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '3.692.93'
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("D40").Value = '''35.71'
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").Value = '0.0₃0766'
$ws.Range("E41").Value = '  +3.20%  '
$ws.Range("D42").Value = '''0.0470'
$ws.Range("E42").Value = '  +8.48%  '
$ws.Range("D43").Value = '''3.42'
$ws.Range("E43").Value = '  +3.38%  '
$ws.Range("D44").Value = '''2.79'
$ws.Range("E44").Value = '  +2.67%  '
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("D47").Value = '''2.81'
$ws.Range("E47").Value = '  +4.44%  '
$ws.Range("E48").Value = '  +2.83%  '
$ws.Range("E49").Value = '  +3.41%  '
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").Value = '''131.97'
$ws.Range("E51").Value = '  +1.03%  '
